$d = $word.ActiveDocument

# Locate the paragraph holding the existing tester's contact details
# (name / student id / email) so the new tester line can be anchored
# immediately after it, without depending on a hard-coded paragraph index.
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*mahirul.islam@mail.mcgill.ca*") {
        $anchorPara = $p
        break
    }
}

if ($anchorPara -eq $null) {
    Write-Output "ERROR: could not locate the anchor (Mahirul Islam) paragraph"
} else {
    # Insert a brand-new, blank paragraph right after the existing tester's
    # line; it naturally inherits that paragraph's (non-spacer) formatting,
    # i.e. just the en-CA/fr-FR language mark with no sz/szCs override.
    $anchorPara.Range.InsertParagraphAfter()

    # Re-fetch the collection; the blank paragraph we just created is now the
    # very next paragraph after the anchor.
    $newPara = $anchorPara.Next()

    # Populate it with the new tester's info as two separate runs (name, then
    # id/email), each explicitly tagged en-CA, by inserting a small OOXML
    # fragment into that paragraph's range.
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
        '<w:p>' + `
        '<w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr>' + `
        '<w:t xml:space="preserve">Tiffany Miller                                           </w:t></w:r>' + `
        '<w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr>' + `
        '<w:t>261113912                                    tiffany.miller@mail.mcgill.ca</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($xml)

    # The paragraph that used to sit right after the anchor (a blank spacer
    # paragraph, sz/szCs 10) has now been pushed one slot further down; flip
    # its language mark from fr-FR to en-CA too.
    $spacerPara = $newPara.Next()
    $spacerPara.Range.LanguageID = "en-CA"

    Write-Output ("Anchor   : [" + $anchorPara.Range.Text + "]")
    Write-Output ("New para : [" + $newPara.Range.Text + "]")
    Write-Output ("Spacer   : [" + $spacerPara.Range.Text + "]")
}
